# The commit swaps the "Office Theme" and "Integral" theme content between
# ppt/theme/theme1.xml (the slide master's theme) and ppt/theme/theme2.xml
# (the notes master's theme): theme1.xml ends up with the default Office
# Theme colors, theme2.xml ends up with the Integral theme colors.
#
# Only the <a:clrScheme> colours actually differ between the two themes
# (fonts/format scheme are already identical), so recreate that swap by
# writing the Office Theme's 12 standard theme colours onto the deck's
# (slide master) theme colour scheme, in MsoThemeColorSchemeIndex order:
# 1 Dark1, 2 Light1, 3 Dark2, 4 Light2, 5 Accent1 .. 10 Accent6,
# 11 Hyperlink, 12 FollowedHyperlink. (RGB values encoded the same way the
# VBA RGB() function does: R + G*256 + B*65536.)

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

$cs.Item(1).RGB  = 0x000000   # dk1      000000
$cs.Item(2).RGB  = 0xFFFFFF   # lt1      FFFFFF
$cs.Item(3).RGB  = 0x6A5444   # dk2      44546A
$cs.Item(4).RGB  = 0xE6E6E7   # lt2      E7E6E6
$cs.Item(5).RGB  = 0xD59B5B   # accent1  5B9BD5
$cs.Item(6).RGB  = 0x317DED   # accent2  ED7D31
$cs.Item(7).RGB  = 0xA5A5A5   # accent3  A5A5A5
$cs.Item(8).RGB  = 0x00C0FF   # accent4  FFC000
$cs.Item(9).RGB  = 0xC47244   # accent5  4472C4
$cs.Item(10).RGB = 0x47AD70   # accent6  70AD47
$cs.Item(11).RGB = 0xC16305   # hlink    0563C1
$cs.Item(12).RGB = 0x724F95   # folHlink 954F72
